# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws4 = $wb.Worksheets.Item("全部类型")

# 展览 sheet (sheet1)
$ws1.Range("F3").Value  = 1177
$ws1.Range("F9").Value  = 3171
$ws1.Range("F10").Value = 712
$ws1.Range("F13").Value = 921
$ws1.Range("F16").Value = 1538
$ws1.Range("F19").Value = 45
$ws1.Range("F20").Value = 1331
$ws1.Range("F21").Value = 444
$ws1.Range("F22").Value = 547
$ws1.Range("F23").Value = 234
$ws1.Range("F24").Value = 7917
$ws1.Range("F25").Value = 9249
$ws1.Range("F29").Value = 114
$ws1.Range("F30").Value = 291

# 全部类型 sheet (sheet4)
$ws4.Range("F4").Value  = 1177
$ws4.Range("F11").Value = 3171
$ws4.Range("F12").Value = 712
$ws4.Range("F15").Value = 921
$ws4.Range("F18").Value = 1538
$ws4.Range("F22").Value = 45
$ws4.Range("F24").Value = 1331
$ws4.Range("F25").Value = 444
$ws4.Range("F26").Value = 547
$ws4.Range("F27").Value = 234
$ws4.Range("F28").Value = 7917
$ws4.Range("F29").Value = 9249
$ws4.Range("F35").Value = 114
$ws4.Range("F36").Value = 291

$wb.Save()
